# frameData.xlsx: change standing right punch animation and adjust move balance
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# "Standing Right Punch" row (row 28): onhit balance value 9 -> 8.
# Dependent formula cells (C4, F4, F28, ...) recalc automatically.
$ws.Range("C28").Value = 8

# Move the active selection/cursor to C29 (matches the author's final cursor position).
$ws.Range("C29").Select()
